$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Range("F4").Value = 3751
$ws1.Range("F5").Value = 3751
$ws1.Range("F6").Value = 291
$ws1.Range("F7").Value = 5286
$ws1.Range("F8").Value = 583
$ws1.Range("F9").Value = 418
$ws1.Range("F10").Value = 226
$ws1.Range("F11").Value = 1052
$ws1.Range("F13").Value = 137
$ws1.Range("F14").Value = 49
$ws1.Range("F16").Value = 360
$ws1.Range("F17").Value = 46
$ws1.Range("F19").Value = 174
$ws1.Range("F22").Value = 6046
$ws1.Range("F23").Value = 6046
$ws1.Range("F26").Value = 14
$ws1.Range("F27").Value = 6806
$ws1.Range("F28").Value = 24
$ws1.Range("F29").Value = 23
$ws1.Range("F30").Value = 3254
$ws1.Range("F31").Value = 366
$ws1.Range("F33").Value = 4458
$ws1.Range("F34").Value = 324
$ws1.Range("F37").Value = 1140
$ws1.Range("F38").Value = 103
$ws1.Range("F41").Value = 920
$ws1.Range("F42").Value = 1124
$ws1.Range("F43").Value = 2057
$ws2.Range("F3").Value = 30
$ws3.Range("F3").Value = 1152
$ws3.Range("F4").Value = 54
$ws4.Range("F4").Value = 1152
$ws4.Range("F5").Value = 54
$ws4.Range("F7").Value = 3751
$ws4.Range("F8").Value = 3751
$ws4.Range("F9").Value = 291
$ws4.Range("F10").Value = 5286
$ws4.Range("F11").Value = 583
$ws4.Range("F12").Value = 418
$ws4.Range("F13").Value = 226
$ws4.Range("F14").Value = 1052
$ws4.Range("F16").Value = 137
$ws4.Range("F17").Value = 49
$ws4.Range("F19").Value = 360
$ws4.Range("F20").Value = 46
$ws4.Range("F23").Value = 174
$ws4.Range("F26").Value = 6046
$ws4.Range("F28").Value = 43
$ws4.Range("F29").Value = 14
$ws4.Range("F30").Value = 6806
$ws4.Range("F31").Value = 24
$ws4.Range("F32").Value = 23
$ws4.Range("F33").Value = 3254
$ws4.Range("F34").Value = 366
$ws4.Range("F36").Value = 4458
$ws4.Range("F37").Value = 324
$ws4.Range("F38").Value = 30
$ws4.Range("F41").Value = 1140
$ws4.Range("F42").Value = 103
$ws4.Range("F45").Value = 920
$ws4.Range("F46").Value = 1124
$ws4.Range("F48").Value = 2057

$wb.Save()
Write-Output "Done: applied 62 changes"